$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.884.61'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '1.866.59'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '305.84'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.5096'
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("D8").Value = '0.3655'
$ws.Range("E8").Value = '  -2.65%  '
$ws.Range("D9").Value = '0.07172'
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("D10").Value = '0.8892'
$ws.Range("E10").Value = '  -0.24%  '
$ws.Range("D11").Value = '20.57'
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.880.29'
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.07484'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '94.35'
$ws.Range("E14").Value = '  +5.28%  '
$ws.Range("D15").Value = '5.218'
$ws.Range("E15").Value = '  -1.58%  '
$ws.Range("E16").Value = '  +0.04%  '
$ws.Range("D17").Value = '0.000008476'
$ws.Range("E17").Value = '  -0.13%  '
$ws.Range("D18").Value = '14.12'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '26.931.21'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.000'
$ws.Range("D21").Style = $ws.Range("D22").Style
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("D22").Value = '2.111.57'
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").Value = '10.32'
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("D24").Value = '6.369'
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").Value = '147.64'
$ws.Range("E25").Value = '  +0.96%  '
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("D27").Value = '17.84'
$ws.Range("E27").Value = '  -0.68%  '
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  +0.43%  '
$ws.Range("D30").Value = '4.671'
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").Value = '4.691'
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("D32").Value = '0.09097'
$ws.Range("E32").Value = '  -1.49%  '
$ws.Range("D33").Value = '0.05026'
$ws.Range("E33").Value = '  -1.93%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7470'
$ws.Range("D34").Style = $ws.Range("D33").Style
$ws.Range("E34").Value = '  +3.01%  '
$ws.Range("D35").Value = '2.983'
$ws.Range("E35").Value = '  -3.17%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = '3.218'
$ws.Range("E37").Value = '  +3.68%  '
$ws.Range("D38").Value = '2.506'
$ws.Range("E38").Value = '  -0.31%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.5573'
$ws.Range("E39").Value = '  +5.45%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01986'
$ws.Range("E40").Value = '  -2.22%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").Value = '6.574'
$ws.Range("E42").Value = '  +1.07%  '
$ws.Range("D43").Value = '115.47'
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").Value = '8.602'
$ws.Range("E44").Value = '  +3.56%  '
$ws.Range("D45").Value = '0.1482'
$ws.Range("E45").Value = '  +0.85%  '
$ws.Range("D46").Value = '0.4753'
$ws.Range("E46").Value = '  +2.96%  '
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").Value = '10.04'
$ws.Range("E48").Value = '  +0.67%  '
$ws.Range("D49").Value = '37.06'
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("D50").Value = '1.552'
$ws.Range("E50").Value = '  -0.40%  '
$ws.Range("D51").Value = '62.91'
$ws.Range("E51").Value = '  -1.14%  '
